$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; D='69.188.23'; E='  +2.06%  '},
    @{Row=3; D='3.826.53'; E='  +0.79%  '},
    @{Row=4; E='  +0.32%  '},
    @{Row=5; D='''627.57'; E='  +4.65%  '},
    @{Row=6; D='''166.05'; E='  +0.46%  '},
    @{Row=7; D='3.823.73'; E='  +0.81%  '},
    @{Row=8; D='''1.00'; E='  -0.09%  '},
    @{Row=9; D='''0.521'; E='  +0.77%  '},
    @{Row=10; E='  +1.77%  '},
    @{Row=11; E='  +0.70%  '},
    @{Row=12; E='  +1.88%  '},
    @{Row=13; E='  +0.90%  '},
    @{Row=14; D='''36.15'; E='  +0.86%  '},
    @{Row=15; D='4.470.55'; E='  +0.78%  '},
    @{Row=16; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='3.813.58'; E='  +0.69%  '},
    @{Row=17; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='69.227.04'; E='  +2.05%  '},
    @{Row=18; D='''18.14'; E='  -1.16%  '},
    @{Row=19; D='''7.15'; E='  +1.23%  '},
    @{Row=20; E='  -0.04%  '},
    @{Row=21; D='''466.84'; E='  +0.85%  '},
    @{Row=22; E='  -1.14%  '},
    @{Row=23; D='''0.709'; E='  +1.41%  '},
    @{Row=24; D='''0.0000153'; E='  +4.74%  '},
    @{Row=25; D='''83.96'; E='  +1.41%  '},
    @{Row=26; D='''12.04'; E='  -0.09%  '},
    @{Row=27; E='  +2.51%  '},
    @{Row=28; E='  +0.34%  '},
    @{Row=29; E='  +0.06%  '},
    @{Row=30; E='  +0.77%  '},
    @{Row=31; E='  -0.27%  '},
    @{Row=32; E='  +0.98%  '},
    @{Row=33; D='''7.32'; E='  -1.65%  '},
    @{Row=34; D='''29.28'; E='  +0.30%  '},
    @{Row=35; D='''9.12'; E='  +1.03%  '},
    @{Row=36; E='  +0.32%  '},
    @{Row=37; E='  +2.28%  '},
    @{Row=38; E='  +7.63%  '},
    @{Row=39; E='  +5.83%  '},
    @{Row=40; D='''5.93'; E='  +2.79%  '},
    @{Row=41; D='''0.984'; E='  -0.20%  '},
    @{Row=42; E='  +0.15%  '},
    @{Row=43; E='  +0.02%  '},
    @{Row=44; D='''1.45'; E='  +5.34%  '},
    @{Row=45; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='''154.85'; E='  +2.37%  '},
    @{Row=46; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='''0.300'; E='  +0.56%  '},
    @{Row=47; D='''46.95'; E='  -1.33%  '},
    @{Row=48; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='''8.48'; E='  +1.53%  '},
    @{Row=49; B='Arweave'; C='https://coinranking.com/coin/7XWg41D1+arweave-ar'; D='''42.55'; E='  -4.85%  '},
    @{Row=50; D='''1.89'; E='  +2.27%  '},
    @{Row=51; D='''0.000283'; E='  +13.94%  '}
)

foreach ($item in $changes) {
    $r = $item.Row
    if ($item.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $item.C }
    if ($item.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $item.D }
    if ($item.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $item.E }
}

"Done"